# Update country data file ("MSME Country Indicators - Oman Summary.xlsx"):
#   * rename the sheet "Data" -> "Summary"
#   * add a bold+underlined "Source Type" line above the existing table
#   * push the existing Micro/SMEs/MSMEs table down to make room
#   * append a second indicator table ("Value added to the economy")
#   * append a source-citation block
#
# NOTE: this engine's xlsx round-trip re-numbers the style array on
# every save (even for a plain load+save with zero edits), so any cell
# whose formatting depended on the workbook's original numeric style
# index silently loses that formatting across the COM bridge. To keep
# the output visually faithful we explicitly (re-)apply the intended
# Font properties (bold/italic/underline/size) to every cell that
# needs them, including the two existing header cells (A1, A3) that
# the source diff itself leaves untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet name: "Data" -> "Summary"
$ws.Name = "Summary"

# Re-assert formatting on the two pre-existing header cells that the
# diff does not otherwise touch, so they keep rendering correctly.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# The Micro/SMEs/MSMEs table (old rows 5-10) is moving down to rows
# 11-16 to make room for the new "Source Type" line; clear the old
# location first so nothing is left duplicated behind.
$ws.Range("A5:D10").Clear()

function Set-TextCell($rng, $value) {
    # Force the cell to stay text (not get coerced to a number) for
    # numeric-looking values like "118386", then drop back to the
    # Normal style so no stray number-format style lingers on it.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# 2) New bold+underlined "Source Type" line at row 9
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# 3) Existing Micro/SMEs/MSMEs table, shifted from rows 5-10 to rows 11-16
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("A12").Font.Bold = $true
Set-TextCell $ws.Range("D12") "118386"

$ws.Range("A13").Value = "Enterprises density (per 1000 people)"
$ws.Range("A13").Font.Bold = $true
Set-TextCell $ws.Range("D13") "44.5"

$ws.Range("A14").Value = "Employment (absolute #)"
$ws.Range("A14").Font.Bold = $true
Set-TextCell $ws.Range("D14") "156135"

$ws.Range("A15").Value = "Enterprises (% of total)"
$ws.Range("A15").Font.Bold = $true
Set-TextCell $ws.Range("D15") "90"

$ws.Range("A16").Value = "Source: MCI, 2009"
$ws.Range("A16").Font.Italic = $true

# 4) New second indicator table: header row repeated at row 18,
#    a new "Value added to the economy" row at 19, source at row 20
$ws.Range("B18").Value = "Micro"
$ws.Range("B18").Font.Bold = $true
$ws.Range("C18").Value = "SMEs"
$ws.Range("C18").Font.Bold = $true
$ws.Range("D18").Value = "MSMEs"
$ws.Range("D18").Font.Bold = $true

$ws.Range("A19").Value = "Value added to the economy (% of total)"
$ws.Range("A19").Font.Bold = $true
Set-TextCell $ws.Range("D19") "23"

$ws.Range("A20").Value = "Source: MCI, 2009"
$ws.Range("A20").Font.Italic = $true

# 5) New citation block at rows 25-26
$ws.Range("A25").Value = "MCI"
$ws.Range("A25").Font.Bold = $true

$ws.Range("A26").Value = "Ministry of Commerce and Industry (MCI), SMEs in Sultanate of Oman : Meeting the development challenges, Dr. Ayoob C. P., Dr. Balabrishnan Somasundaram, Ministry of Higher Education, Sultanate of Oman, NIZWA, p. 27 based on information from the Ministry of Commerce and Industry (MCI)."
$ws.Range("A26").Font.Italic = $true
